$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (4) = "algorithm", Column S (19) = "runcommand".
# Rows 18-33 invoke main_program.py --mip_emphasis 1 ...
# Rows 34-49 invoke main_program_one_depth_cascade.py ...
# Both need "--timelimit 12 --export_results_file" inserted right after the
# script-selection tokens ("--mip_emphasis 1" or the script name itself).

for ($r = 18; $r -le 33; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $s = $ws.Cells.Item($r, 19).Value2

    $newD = $d.Replace("main_program.py --mip_emphasis 1", "main_program.py --mip_emphasis 1 --timelimit 12 --export_results_file")
    $newS = $s.Replace("main_program.py --mip_emphasis 1", "main_program.py --mip_emphasis 1 --timelimit 12 --export_results_file")

    $ws.Cells.Item($r, 4).Value = $newD
    $ws.Cells.Item($r, 19).Value = $newS
}

for ($r = 34; $r -le 49; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $s = $ws.Cells.Item($r, 19).Value2

    $newD = $d.Replace("main_program_one_depth_cascade.py", "main_program_one_depth_cascade.py --timelimit 12 --export_results_file")
    $newS = $s.Replace("main_program_one_depth_cascade.py", "main_program_one_depth_cascade.py --timelimit 12 --export_results_file")

    $ws.Cells.Item($r, 4).Value = $newD
    $ws.Cells.Item($r, 19).Value = $newS
}
